$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing data row (old row 2 -> row 3)
$ws.Rows("2:2").Insert()

# Force NDC / Lot / Expiry columns to text so numeric/date-looking values
# are stored as strings (matches the inlineStr cells in the target file)
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

# New record inserted at row 2
$ws.Range("A2").Value = "6787743305"
$ws.Range("B2").Value = "Ascend Laboratories, LLC"
$ws.Range("C2").Value = "Aripiprazole"
$ws.Range("D2").Value = "15 mg/1"
$ws.Range("E2").Value = "22140477"
$ws.Range("F2").Value = "01/24/31"
$ws.Range("G2").Value = 5

# Existing record (now on row 3) had its QTY updated from 1 to 3
$ws.Range("G3").Value = 3

# Re-materialize the trailing blank row (was row 3, now row 4) so the
# sheet's used range / dimension extends to row 4, matching the captured
# full-data-state snapshot. Re-assigning a property to its own current
# (no-op) value touches the cell/row without introducing new formatting.
$fontSize = $ws.Cells.Item(4, 1).Font.Size
$ws.Cells.Item(4, 1).Font.Size = $fontSize
